$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as text in the original workbook (t="inlineStr"),
# so force text format before assigning numeric-looking strings to avoid numeric coercion.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "264.00"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.72"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.213"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06087"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.725"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8153"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1590"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08144"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03152"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09269"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.927"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04829"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006243"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006199"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006113"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001099"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001504"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.267"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3393"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1270"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002687"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04640"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007249"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003138"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01050"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006100"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000752"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7516"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03253"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002104"
$ws.Range("E49").Value = "48CryptobidCoinCBCBestin24h"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01243"
